# Update Argent (silver) / solar prices for 2025-03-05 (row 4, column B)
# across the relevant sheets. Each value is stored as text (it may look
# like a number, e.g. "5,211" or "0.293"), so we force the cell's number
# format to Text ("@") before assigning the value. This prevents Excel
# from re-interpreting the text as a numeric value (which would also
# silently rewrite the other, untouched cells' displayed formatting).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cell Topcon 183mm")
$rng = $ws.Range("B4")
$rng.NumberFormat = "@"
$rng.Value = "0.293"

$ws = $wb.Worksheets.Item("Silver Rear_side")
$rng = $ws.Range("B4")
$rng.NumberFormat = "@"
$rng.Value = "5,211"

$ws = $wb.Worksheets.Item("Silver Busbar front-side")
$rng = $ws.Range("B4")
$rng.NumberFormat = "@"
$rng.Value = "7,801"

$ws = $wb.Worksheets.Item("Silver finger front-side")
$rng = $ws.Range("B4")
$rng.NumberFormat = "@"
$rng.Value = "7,851"

$ws = $wb.Worksheets.Item("USD_CNY")
$rng = $ws.Range("B4")
$rng.NumberFormat = "@"
$rng.Value = "7.2842"
